$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.161.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.001.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +13.15%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.998.54"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.45%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.182.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.501.44"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.005.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.13"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.55%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.90%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +15.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.23%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +17.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +20.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000106"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.60"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.17"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.20%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.78%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.993"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.16"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +14.85%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.91"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +17.50%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.89"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.40"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.29"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +14.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.791.59"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.65"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.83"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +11.98%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.83%  "
